$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2578893899917603
$ws.Range("B1").Value = 0.1907126903533936
$ws.Range("C1").Value = 0.299547404050827
$ws.Range("D1").Value = 3.94853949546814
$ws.Range("E1").Value = 1.696321964263916
